$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(403, 366),
    @(807, 70),
    @(616, 84),
    @(674, 456),
    @(580, 318),
    @(145, 302),
    @(964, 322),
    @(477, 447),
    @(1002, 116),
    @(151, 214),
    @(782, 199),
    @(263, 88),
    @(201, 6),
    @(459, 158),
    @(838, 288),
    @(875, 153),
    @(589, 462),
    @(699, 14),
    @(520, 211),
    @(298, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
